$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")

# ALC row 112
$ws_ALC.Range("H112").Value = 1172.4
$ws_ALC.Range("I112").Value = 200
$ws_ALC.Range("J112").Value = 1415.5
$ws_ALC.Range("K112").Value = 600
$ws_ALC.Range("L112").Value = 4246.5
$ws_ALC.Range("M112").Value = 508
$ws_ALC.Range("N112").Value = -6462.5

# ALC row 113
$ws_ALC.Range("H113").Value = 2874.5652
$ws_ALC.Range("J113").Value = 2888.6667
$ws_ALC.Range("L113").Value = 2888.6667
$ws_ALC.Range("N113").Value = -9396.6667

# ARM row 2
$ws_ARM.Range("H2").Value = 602.5
$ws_ARM.Range("I2").Value = 637
$ws_ARM.Range("J2").Value = 499
$ws_ARM.Range("K2").Value = 637
$ws_ARM.Range("L2").Value = 499
$ws_ARM.Range("M2").Value = -524
$ws_ARM.Range("N2").Value = -725

# ARM row 74
$ws_ARM.Range("H74").Value = 5104127.5
$ws_ARM.Range("I74").Value = 10002271
$ws_ARM.Range("J74").Value = 1894.2916
$ws_ARM.Range("K74").Value = 10002271
$ws_ARM.Range("L74").Value = 1894.2916
$ws_ARM.Range("M74").Value = -10001397
$ws_ARM.Range("N74").Value = -3642.2916

# ARM row 77
$ws_ARM.Range("H77").Value = 5104127.5
$ws_ARM.Range("I77").Value = 10002271
$ws_ARM.Range("J77").Value = 1894.2916
$ws_ARM.Range("K77").Value = 50011355
$ws_ARM.Range("L77").Value = 9471.458000000001
$ws_ARM.Range("M77").Value = -50006987
$ws_ARM.Range("N77").Value = -18207.458

# ARM row 110
$ws_ARM.Range("H110").Value = 7523.0454
$ws_ARM.Range("I110").Value = 9208
$ws_ARM.Range("K110").Value = 9208
$ws_ARM.Range("M110").Value = -7163

# ARM row 116
$ws_ARM.Range("H116").Value = 602.5
$ws_ARM.Range("I116").Value = 637
$ws_ARM.Range("J116").Value = 499
$ws_ARM.Range("K116").Value = 637
$ws_ARM.Range("L116").Value = 499
$ws_ARM.Range("M116").Value = 1657
$ws_ARM.Range("N116").Value = -5087

# ARM row 122
$ws_ARM.Range("H122").Value = 2831.6
$ws_ARM.Range("I122").Value = 2452.182
$ws_ARM.Range("J122").Value = 3875
$ws_ARM.Range("K122").Value = 7356.545999999999
$ws_ARM.Range("L122").Value = 11625
$ws_ARM.Range("M122").Value = -4906.545999999999
$ws_ARM.Range("N122").Value = -16525

# BSM row 3
$ws_BSM.Range("H3").Value = 602.5
$ws_BSM.Range("I3").Value = 637
$ws_BSM.Range("J3").Value = 499
$ws_BSM.Range("K3").Value = 637
$ws_BSM.Range("L3").Value = 499
$ws_BSM.Range("M3").Value = -523
$ws_BSM.Range("N3").Value = -727

# BSM row 86
$ws_BSM.Range("H86").Value = 20002140
$ws_BSM.Range("I86").Value = 28573184
$ws_BSM.Range("J86").Value = 3035.6667
$ws_BSM.Range("K86").Value = 28573184
$ws_BSM.Range("L86").Value = 3035.6667
$ws_BSM.Range("M86").Value = -28572061
$ws_BSM.Range("N86").Value = -5281.6667

# BSM row 89
$ws_BSM.Range("H89").Value = 20002140
$ws_BSM.Range("I89").Value = 28573184
$ws_BSM.Range("J89").Value = 3035.6667
$ws_BSM.Range("K89").Value = 142865920
$ws_BSM.Range("L89").Value = 15178.3335
$ws_BSM.Range("M89").Value = -142860304
$ws_BSM.Range("N89").Value = -26410.3335

# CRP row 31
$ws_CRP.Range("H31").Value = 3527823.2
$ws_CRP.Range("I31").Value = 5743802
$ws_CRP.Range("J31").Value = 2402.5
$ws_CRP.Range("K31").Value = 5743802
$ws_CRP.Range("L31").Value = 2402.5
$ws_CRP.Range("M31").Value = -5743507
$ws_CRP.Range("N31").Value = -2992.5

# CRP row 34
$ws_CRP.Range("H34").Value = 3527823.2
$ws_CRP.Range("I34").Value = 5743802
$ws_CRP.Range("J34").Value = 2402.5
$ws_CRP.Range("K34").Value = 5743802
$ws_CRP.Range("L34").Value = 2402.5
$ws_CRP.Range("M34").Value = -5743600
$ws_CRP.Range("N34").Value = -2806.5

# CRP row 58
$ws_CRP.Range("H58").Value = 2093.2778
$ws_CRP.Range("I58").Value = 1344.25
$ws_CRP.Range("J58").Value = 2307.2856
$ws_CRP.Range("K58").Value = 1344.25
$ws_CRP.Range("L58").Value = 2307.2856
$ws_CRP.Range("M58").Value = -1141.25
$ws_CRP.Range("N58").Value = -2713.2856

# CRP row 63
$ws_CRP.Range("H63").Value = 15500
$ws_CRP.Range("J63").Value = 15500
$ws_CRP.Range("L63").Value = 15500
$ws_CRP.Range("N63").Value = -16872

# CRP row 66
$ws_CRP.Range("H66").Value = 15500
$ws_CRP.Range("J66").Value = 15500
$ws_CRP.Range("L66").Value = 46500
$ws_CRP.Range("N66").Value = -53364

# CRP row 122
$ws_CRP.Range("H122").Value = 807.6875
$ws_CRP.Range("I122").Value = 794.9
$ws_CRP.Range("J122").Value = 829
$ws_CRP.Range("K122").Value = 2384.7
$ws_CRP.Range("L122").Value = 2487
$ws_CRP.Range("M122").Value = 65.30000000000018
$ws_CRP.Range("N122").Value = -7387

# CRP row 132
$ws_CRP.Range("H132").Value = 3080.8
$ws_CRP.Range("I132").Value = 2842.5
$ws_CRP.Range("J132").Value = 3353.1428
$ws_CRP.Range("K132").Value = 8527.5
$ws_CRP.Range("L132").Value = 10059.4284
$ws_CRP.Range("M132").Value = -5997.5
$ws_CRP.Range("N132").Value = -15119.4284

# CRP row 134
$ws_CRP.Range("H134").Value = 2931.9
$ws_CRP.Range("I134").Value = 3255.318
$ws_CRP.Range("J134").Value = 2042.5
$ws_CRP.Range("K134").Value = 9765.954000000002
$ws_CRP.Range("L134").Value = 6127.5
$ws_CRP.Range("M134").Value = -7230.954000000002
$ws_CRP.Range("N134").Value = -11197.5

# CRP row 136
$ws_CRP.Range("H136").Value = 2093.2778
$ws_CRP.Range("I136").Value = 1344.25
$ws_CRP.Range("J136").Value = 2307.2856
$ws_CRP.Range("K136").Value = 4032.75
$ws_CRP.Range("L136").Value = 6921.8568
$ws_CRP.Range("M136").Value = -1482.75
$ws_CRP.Range("N136").Value = -12021.8568

# CUL row 105
$ws_CUL.Range("H105").Value = 35731350
$ws_CUL.Range("I105").Value = 0
$ws_CUL.Range("J105").Value = 35731350
$ws_CUL.Range("K105").Value = 0
$ws_CUL.Range("L105").Value = 107194050
$ws_CUL.Range("M105").ClearContents()
$ws_CUL.Range("N105").Value = -107199292

# CUL row 114
$ws_CUL.Range("H114").Value = 1206.4348
$ws_CUL.Range("J114").Value = 608.3333
$ws_CUL.Range("L114").Value = 1824.9999
$ws_CUL.Range("N114").Value = -8332.999900000001

# CUL row 121
$ws_CUL.Range("H121").Value = 2156.647
$ws_CUL.Range("I121").Value = 496
$ws_CUL.Range("J121").Value = 2848.5833
$ws_CUL.Range("K121").Value = 1488
$ws_CUL.Range("L121").Value = 8545.749899999999
$ws_CUL.Range("M121").Value = -178
$ws_CUL.Range("N121").Value = -11165.7499

# CUL row 132
$ws_CUL.Range("H132").Value = 1200
$ws_CUL.Range("I132").Value = 700
$ws_CUL.Range("J132").Value = 1950
$ws_CUL.Range("K132").Value = 6300
$ws_CUL.Range("L132").Value = 17550
$ws_CUL.Range("M132").Value = -3770
$ws_CUL.Range("N132").Value = -22610

# CUL row 137
$ws_CUL.Range("H137").Value = 8388.261
$ws_CUL.Range("I137").Value = 1200
$ws_CUL.Range("J137").Value = 8715
$ws_CUL.Range("K137").Value = 3600
$ws_CUL.Range("L137").Value = 26145
$ws_CUL.Range("M137").Value = 1500
$ws_CUL.Range("N137").Value = -36345

# GSM row 102
$ws_GSM.Range("H102").Value = 2528.75
$ws_GSM.Range("I102").Value = 2806.4285
$ws_GSM.Range("J102").Value = 2140
$ws_GSM.Range("K102").Value = 2806.4285
$ws_GSM.Range("L102").Value = 2140
$ws_GSM.Range("M102").Value = -1184.4285
$ws_GSM.Range("N102").Value = -5384

# GSM row 122
$ws_GSM.Range("H122").Value = 3972.0356
$ws_GSM.Range("I122").Value = 3421.158
$ws_GSM.Range("J122").Value = 5135
$ws_GSM.Range("K122").Value = 10263.474
$ws_GSM.Range("L122").Value = 15405
$ws_GSM.Range("M122").Value = -7813.474
$ws_GSM.Range("N122").Value = -20305
